$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D stays text (avoid Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '23.368.24'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '1.631.48'
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '0.9973'
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").Value = '303.16'
$ws.Range("E6").Value = '  -0.73%  '

$ws.Range("D7").Value = '0.3752'

$ws.Range("D8").Value = '51.67'
$ws.Range("E8").Value = '  -0.40%  '

$ws.Range("D9").Value = '0.3569'

$ws.Range("D10").Value = '0.08176'
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("D11").Value = '1.231'
$ws.Range("E11").Value = '  -0.06%  '

$ws.Range("D12").Value = '0.9988'
$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("D13").Value = '22.26'
$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").Value = '6.487'
$ws.Range("E14").Value = '  -1.50%  '

$ws.Range("D15").Value = '7.320'
$ws.Range("E15").Value = '  +0.41%  '

$ws.Range("D16").Value = '0.00001222'
$ws.Range("E16").Value = '  -2.38%  '

$ws.Range("D17").Value = '1.622.11'
$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("D18").Value = '95.54'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").Value = '0.06938'
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").Value = '6.699'
$ws.Range("E20").Value = '  +3.42%  '

$ws.Range("D21").Value = '17.43'
$ws.Range("E21").Value = '  -2.13%  '

$ws.Range("D22").Value = '0.9978'
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("D23").Value = '12.43'
$ws.Range("E23").Value = '  -2.79%  '

$ws.Range("D24").Value = '23.352.64'
$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("D25").Value = '2.518'
$ws.Range("E25").Value = '  +2.19%  '

$ws.Range("D26").Value = '3.089'
$ws.Range("E26").Value = '  -3.35%  '

$ws.Range("D27").Value = '21.09'
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("D28").Value = '153.13'
$ws.Range("E28").Value = '  +2.13%  '

$ws.Range("D29").Value = '5.177'
$ws.Range("E29").Value = '  -2.88%  '

$ws.Range("D30").Value = '133.80'
$ws.Range("E30").Value = '  -0.85%  '

$ws.Range("D31").Value = '1.801.00'
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").Value = '1.095'
$ws.Range("E32").Value = '  +6.57%  '

$ws.Range("D33").Value = '6.620'
$ws.Range("E33").Value = '  -3.24%  '

$ws.Range("D34").Value = '11.60'
$ws.Range("E34").Value = '  +5.93%  '

$ws.Range("D35").Value = '2.022'
$ws.Range("E35").Value = '  -11.08%  '

$ws.Range("D36").Value = '0.02739'
$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("D37").Value = '0.08755'
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").Value = '0.2482'
$ws.Range("E38").Value = '  -2.24%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.936'
$ws.Range("E39").Value = '  -2.47%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06949'
$ws.Range("E40").Value = '  -2.74%  '

$ws.Range("D41").Value = '12.54'
$ws.Range("E41").Value = '  +1.55%  '

$ws.Range("D42").Value = '0.6946'
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("D43").Value = '1.327'
$ws.Range("E43").Value = '  -2.31%  '

$ws.Range("D44").Value = '15.58'
$ws.Range("E44").Value = '  -4.30%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.6396'
$ws.Range("E45").Value = '  -2.54%  '

$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = '0.9975'
$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("D47").Value = '2.291'
$ws.Range("E47").Value = '  -2.12%  '

$ws.Range("D48").Value = '3.957'
$ws.Range("E48").Value = '  -1.10%  '

$ws.Range("D49").Value = '0.07919'
$ws.Range("E49").Value = '  -1.60%  '

$ws.Range("D50").Value = '127.56'
$ws.Range("E50").Value = '  +1.57%  '

$ws.Range("D51").Value = '1.176'
$ws.Range("E51").Value = '  -2.94%  '

# Reset style on column D back to default (Normal) now that values are set as text
$ws.Range("D2:D51").Style = "Normal"
